# Auto-generated update of computed market-price / profit columns (H:N)
# across all item sheets, per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 309.46155
$ws.Range("I2").Value = 279.2
$ws.Range("J2").Value = 328.375
$ws.Range("K2").Value = 279.2
$ws.Range("L2").Value = 328.375
$ws.Range("M2").Value = -166.2
$ws.Range("N2").Value = -554.375
$ws.Range("H17").Value = 1460.8088
$ws.Range("J17").Value = 1245.591
$ws.Range("L17").Value = 3736.773
$ws.Range("N17").Value = -4072.773
$ws.Range("H28").Value = 474.25
$ws.Range("I28").Value = 348.75
$ws.Range("K28").Value = 348.75
$ws.Range("M28").Value = 136.25
$ws.Range("H51").Value = 2643
$ws.Range("I51").Value = 1225
$ws.Range("J51").Value = 2997.5
$ws.Range("K51").Value = 1225
$ws.Range("L51").Value = 2997.5
$ws.Range("M51").Value = -741
$ws.Range("N51").Value = -3965.5
$ws.Range("H94").Value = 2999.25
$ws.Range("I94").Value = 2999.25
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 2999.25
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -2548.25
$ws.Range("N94").ClearContents()
$ws.Range("H98").Value = 1861.6364
$ws.Range("I98").Value = 1947.5555
$ws.Range("K98").Value = 1947.5555
$ws.Range("M98").Value = -449.5554999999999
$ws.Range("H106").Value = 1665.875
$ws.Range("I106").Value = 1665.875
$ws.Range("K106").Value = 1665.875
$ws.Range("M106").Value = -1034.875
$ws.Range("H113").Value = 51176.25
$ws.Range("J113").Value = 2352.5
$ws.Range("L113").Value = 2352.5
$ws.Range("N113").Value = -8860.5
$ws.Range("H121").Value = 1355.5
$ws.Range("J121").Value = 1596.6
$ws.Range("L121").Value = 4789.799999999999
$ws.Range("N121").Value = -8283.799999999999
$ws.Range("H122").Value = 1861.6364
$ws.Range("I122").Value = 1947.5555
$ws.Range("K122").Value = 5842.666499999999
$ws.Range("M122").Value = -3392.666499999999
$ws.Range("H132").Value = 1156.0526
$ws.Range("I132").Value = 1069.0577
$ws.Range("J132").Value = 2060.8
$ws.Range("K132").Value = 3207.1731
$ws.Range("L132").Value = 6182.400000000001
$ws.Range("M132").Value = -677.1731
$ws.Range("N132").Value = -11242.4
$ws.Range("H138").Value = 3894.652
$ws.Range("I138").Value = 4583.923
$ws.Range("J138").Value = 2998.6
$ws.Range("K138").Value = 13751.769
$ws.Range("L138").Value = 8995.799999999999
$ws.Range("M138").Value = -8611.769
$ws.Range("N138").Value = -19275.8

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2872.4033
$ws.Range("I32").Value = 2579.2354
$ws.Range("K32").Value = 2579.2354
$ws.Range("M32").Value = -2292.2354
$ws.Range("H74").Value = 1218.9688
$ws.Range("I74").Value = 941
$ws.Range("J74").Value = 1929.3334
$ws.Range("K74").Value = 941
$ws.Range("L74").Value = 1929.3334
$ws.Range("M74").Value = -67
$ws.Range("N74").Value = -3677.3334
$ws.Range("H77").Value = 1218.9688
$ws.Range("I77").Value = 941
$ws.Range("J77").Value = 1929.3334
$ws.Range("K77").Value = 4705
$ws.Range("L77").Value = 9646.666999999999
$ws.Range("M77").Value = -337
$ws.Range("N77").Value = -18382.667
$ws.Range("H122").Value = 1966.8125
$ws.Range("I122").Value = 1962.1428
$ws.Range("J122").Value = 1999.5
$ws.Range("K122").Value = 5886.428400000001
$ws.Range("L122").Value = 5998.5
$ws.Range("M122").Value = -3436.428400000001
$ws.Range("N122").Value = -10898.5
$ws.Range("H132").Value = 1582.7441
$ws.Range("I132").Value = 1078.9395
$ws.Range("J132").Value = 3245.3
$ws.Range("K132").Value = 3236.8185
$ws.Range("L132").Value = 9735.900000000001
$ws.Range("M132").Value = -706.8184999999999
$ws.Range("N132").Value = -14795.9

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 136386.53
$ws.Range("I86").Value = 3579.8
$ws.Range("J86").Value = 402000
$ws.Range("K86").Value = 3579.8
$ws.Range("L86").Value = 402000
$ws.Range("M86").Value = -2456.8
$ws.Range("N86").Value = -404246
$ws.Range("H89").Value = 136386.53
$ws.Range("I89").Value = 3579.8
$ws.Range("J89").Value = 402000
$ws.Range("K89").Value = 17899
$ws.Range("L89").Value = 2010000
$ws.Range("M89").Value = -12283
$ws.Range("N89").Value = -2021232
$ws.Range("H94").Value = 664.8
$ws.Range("I94").Value = 706.8570999999999
$ws.Range("K94").Value = 706.8570999999999
$ws.Range("M94").Value = -255.8570999999999
$ws.Range("H105").Value = 2387.4783
$ws.Range("I105").Value = 2295.65
$ws.Range("K105").Value = 2295.65
$ws.Range("M105").Value = -548.6500000000001
$ws.Range("H107").Value = 3500
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("H134").Value = 4711.8
$ws.Range("I134").Value = 5619.9688
$ws.Range("J134").Value = 2476.3076
$ws.Range("K134").Value = 16859.9064
$ws.Range("L134").Value = 7428.9228
$ws.Range("M134").Value = -14324.9064
$ws.Range("N134").Value = -12498.9228

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2116.158
$ws.Range("I31").Value = 1925.125
$ws.Range("J31").Value = 2255.0908
$ws.Range("K31").Value = 1925.125
$ws.Range("L31").Value = 2255.0908
$ws.Range("M31").Value = -1630.125
$ws.Range("N31").Value = -2845.0908
$ws.Range("H34").Value = 2116.158
$ws.Range("I34").Value = 1925.125
$ws.Range("J34").Value = 2255.0908
$ws.Range("K34").Value = 1925.125
$ws.Range("L34").Value = 2255.0908
$ws.Range("M34").Value = -1723.125
$ws.Range("N34").Value = -2659.0908
$ws.Range("I93").Value = 4950
$ws.Range("K93").Value = 4950
$ws.Range("M93").Value = -3078
$ws.Range("H132").Value = 1470.0588
$ws.Range("I132").Value = 963.5185
$ws.Range("K132").Value = 2890.5555
$ws.Range("M132").Value = -360.5554999999999
$ws.Range("H134").Value = 1486.8334
$ws.Range("I134").Value = 1247.4324
$ws.Range("K134").Value = 3742.2972
$ws.Range("M134").Value = -1207.2972

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1714528.4
$ws.Range("I4").Value = 2500149.8
$ws.Range("J4").Value = 667033.3
$ws.Range("K4").Value = 7500449.399999999
$ws.Range("L4").Value = 2001099.9
$ws.Range("M4").Value = -7500337.399999999
$ws.Range("N4").Value = -2001323.9
$ws.Range("H107").Value = 525.6667
$ws.Range("I107").Value = 361.57144
$ws.Range("K107").Value = 1084.71432
$ws.Range("M107").Value = 835.28568
$ws.Range("H131").Value = 14038.935
$ws.Range("J131").Value = 14483.813
$ws.Range("L131").Value = 43451.439
$ws.Range("N131").Value = -53531.439
$ws.Range("H132").Value = 825
$ws.Range("I132").Value = 641.125
$ws.Range("J132").Value = 1119.2
$ws.Range("K132").Value = 5770.125
$ws.Range("L132").Value = 10072.8
$ws.Range("M132").Value = -3240.125
$ws.Range("N132").Value = -15132.8

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 14725
$ws.Range("J26").Value = 14966.667
$ws.Range("L26").Value = 14966.667
$ws.Range("N26").Value = -15526.667
$ws.Range("H50").Value = 14725
$ws.Range("J50").Value = 14966.667
$ws.Range("L50").Value = 14966.667
$ws.Range("N50").Value = -15962.667
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H122").Value = 1408.5
$ws.Range("I122").Value = 1211.3334
$ws.Range("K122").Value = 3634.0002
$ws.Range("M122").Value = -1184.0002
$ws.Range("H126").Value = 3144437
$ws.Range("I126").Value = 13892242
$ws.Range("J126").Value = 73635.71000000001
$ws.Range("K126").Value = 41676726
$ws.Range("L126").Value = 220907.13
$ws.Range("M126").Value = -41674256
$ws.Range("N126").Value = -225847.13
$ws.Range("H132").Value = 1204585.6
$ws.Range("I132").Value = 1924915.5
$ws.Range("K132").Value = 5774746.5
$ws.Range("M132").Value = -5772216.5
$ws.Range("H133").Value = 85666.664
$ws.Range("J133").Value = 85666.664
$ws.Range("L133").Value = 85666.664
$ws.Range("N133").Value = -95786.664

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 11360.5
$ws.Range("I40").Value = 14644.25
$ws.Range("J40").Value = 6982.1665
$ws.Range("K40").Value = 14644.25
$ws.Range("L40").Value = 6982.1665
$ws.Range("M40").Value = -14508.25
$ws.Range("N40").Value = -7254.1665
$ws.Range("H100").Value = 1516.8334
$ws.Range("I100").Value = 1160.2
$ws.Range("K100").Value = 1160.2
$ws.Range("M100").Value = -619.2
$ws.Range("H122").Value = 8033
$ws.Range("I122").Value = 7246.2
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 21738.6
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -19288.6
$ws.Range("N122").Value = -34900
$ws.Range("H132").Value = 1525.7354
$ws.Range("I132").Value = 1281.36
$ws.Range("J132").Value = 2204.5557
$ws.Range("K132").Value = 3844.08
$ws.Range("L132").Value = 6613.6671
$ws.Range("M132").Value = -1314.08
$ws.Range("N132").Value = -11673.6671
$ws.Range("H136").Value = 2839.5483
$ws.Range("I136").Value = 1769
$ws.Range("J136").Value = 5087.7
$ws.Range("K136").Value = 5307
$ws.Range("L136").Value = 15263.1
$ws.Range("M136").Value = -2757
$ws.Range("N136").Value = -20363.1
$ws.Range("H141").Value = 48713.57
$ws.Range("J141").Value = 48713.57
$ws.Range("L141").Value = 48713.57
$ws.Range("N141").Value = -59073.57

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 71847.37
$ws.Range("I122").Value = 78925.10000000001
$ws.Range("J122").Value = 1070
$ws.Range("K122").Value = 236775.3
$ws.Range("L122").Value = 3210
$ws.Range("M122").Value = -234325.3
$ws.Range("N122").Value = -8110
$ws.Range("H132").Value = 1540.9286
$ws.Range("I132").Value = 860.4483
$ws.Range("K132").Value = 2581.3449
$ws.Range("M132").Value = -51.34490000000005
$ws.Range("H136").Value = 55559896
$ws.Range("I136").Value = 79369570
$ws.Range("K136").Value = 238108710
$ws.Range("M136").Value = -238106160
